# Adding meteogram upgrade script.
# Sort the station table (A1:F22, header in row 1) descending by the
# "Domain" column (C), add the AutoFilter on the header row, and move
# the active selection to F9 - matching the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sort A1:F22 (header row 1) descending by column C ("Domain") ----------
$sortRange = $ws.Range("A1:F22")
$keyRange  = $ws.Range("C1:C22")

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($keyRange, 0, 2)
$sortObj.SetRange($sortRange)
$sortObj.Header = 1
$sortObj.Apply()

# --- AutoFilter on the header row only (A1:F1) ------------------------------
$ws.Range("A1:F1").AutoFilter()

# --- Register the (hidden, sheet-scoped) _FilterDatabase defined name ------
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$F`$1")
$fdName.Visible = $false

# --- Move the active selection to F9 ---------------------------------------
$ws.Range("F9").Select()
